# Releasenote.docx edit:
#   - wrap the existing "Releasenote" heading run in proofErr spellStart/spellEnd markers
#   - append the MVP usage documentation as new body paragraphs after the heading
#
# InsertXML() replaces the contents of the exact Range it is called on, so we target
# the heading paragraph's own Range and replay that paragraph (with its original
# w14:paraId/w:rsid* attributes preserved) immediately followed by the six new
# paragraphs (two of them intentionally empty) that make up the release notes text.

$d = $word.ActiveDocument

$xmlPayload = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4D8586A7" w14:textId="3660F58D" w:rsidR="004135F5" w:rsidRDefault="004135F5" w:rsidP="004135F5"><w:pPr><w:pStyle w:val="berschrift1"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Releasenote</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Zum Starten des MVP muss die .exe-Datei ausgeführt werden. Anschließend sollte sich die Oberfläche der Anwendung öffnen. Hier besteht nun die Möglichkeit Eingaben zu tätigen. Eingaben müssen durch den oben mittig-links platzierten „Eingabeschlitz“ erfolgen. Bei jeder Änderung der Eingabe wird überprüft, ob die Eingabe in einen Kontakt übersetzt werden kann. Sollte dies nicht der Fall sein, so erscheint rechts neben dem Knopf „Hinzufügen“ ein entsprechender Hinweis. Außerdem kann dieser Knopf in diesem Fall nicht verwendet werden. Wenn eine Aufteilung in einen Kontakt jedoch möglich ist, so werden die einzelnen Bestandteile in der Vorschau angezeigt. Die Vorschau befindet sich dabei mittig-rechts. In der Vorschaut besteht die Möglichkeit der Bearbeitung des Vor- und Nachnamens. Außerdem kann direkt die Breifanschrift eingesehen werden. Diese kann zwar bearbeitet werden, wird jedoch später nicht mit abgespeichert. </w:t></w:r><w:r><w:t xml:space="preserve">Wenn eine Eingabe in einen Kontakt übersetzt werden kann, so kann der Knopf „Hinzufügen“ betätigt werden. Sollten vor drücken des Knopfes Änderungen an dem Vor- oder Nachnamen erfolgt sein, so werden diese Änderungen mit abgespeichert. In diesem MVP werden die Daten nicht persistent abgespeichert. Datensätze werden nach dem „hinzufügen“ in der farblich hervorgehobenen Liste angezeigt und sind dort auswählbar. Bei Auswahl erscheinen die Einzelteile des Kontakts erneut in der Vorschau. Hierbei ist darauf zu achten, dass die Oberfläche eine Bearbeitung anbietet, diese jedoch keine Auswirkung auf den abgespeicherten Kontakt hat. Kontakte können ausschließlich hinzugefügt oder gelöscht werden. Jedoch nicht nachträglich bearbeitet. Gelöscht werden kann ein Kontakt, wenn dieser in der Liste ausgewählt und anschließend der Knopf „Löschen“ betätigt wird. Sollte kein Kontakt ausgewählt sein, so kann der Knopf nicht betätigt werden. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Anreden und Titel können durch die entsprechenden Knöpfe unten-links bearbeitet, gelöscht oder hinzugefügt werden. Diese werden hierbei persistent gespeichert (in einer Datei). Die entsprechende Datei befindet sich hierbei im Ordner der .exe-Datei. </w:t></w:r><w:r><w:t xml:space="preserve">Durch Auswahl eines Eintrags der Titel oder Anreden kann dieser Bearbeitet oder gelöscht werden. Bei Betätigung des Knopfes </w:t></w:r><w:r><w:t xml:space="preserve">„Hinzufügen oder Updaten“ wird überprüft, ob ein entsprechender Eintrag schon vorhanden ist. Falls ja, wird dieser aktualisiert. Falls nein, wird ein entsprechender neuer Eintrag hinzugefügt. Durch Betätigung der Schaltfläche „Löschen“, wird ein Eintrag aus der Datei gelöscht (ACHTUNG: nach Löschen eines Eintrags wird dieser nicht mehr erkannt!). </w:t></w:r><w:r><w:t>Erkannt werden Titel oder Anreden nur, wenn diese exakt mit den hier definierten Anreden/Titeln übereinstimmen. Soll ein neuer Titel erkannt werden, so muss dieser vor Eingabe in dem „Eingabeschlitz“ durch die entsprechende Funktionalität hinzugefügt werden.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Die Benutzung wurde anhand dieser „Use Cases“ </w:t></w:r><w:r><w:t xml:space="preserve">durch das Entwicklerteam getestet. Sollte bei anderer Benutzung oder ähnlichem ein Fehlverhalten auftreten, so können Sie sich gerne an das Entwicklungsteam wenden. </w:t></w:r><w:r><w:t>Zu beachten ist hierbei, dass es sich um einen MVP handelt. Dieser stellt keinen Anspruch auf absolute Fehlerfreiheit</w:t></w:r><w:r><w:t xml:space="preserve"> (auch wenn diese </w:t></w:r><w:r><w:t>dennoch</w:t></w:r><w:r><w:t xml:space="preserve"> vom Entwicklerteam gewünscht wäre).</w:t></w:r></w:p>
'@

$headingRange = $d.Paragraphs(1).Range
$headingRange.InsertXML($xmlPayload)
